# "changes in concise marksheet" - update Corr/total marks on the quiz marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> Right count
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right count and the Correct/Total summary string
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
